$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.831.96'
$ws.Range('D3').Value = '1.561.14'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.82'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.486'
$ws.Range('E6').Value = '  -1.30%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.76'
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.246'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0865'
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').Value = '1.785.09'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').Value = '1.565.68'
$ws.Range('E13').Value = '  +0.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.72'
$ws.Range('E14').Value = '  -1.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.513'
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').Value = '26.851.21'
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.21'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '214.70'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.34'
$ws.Range('E19').Value = '  +1.81%  '
$ws.Range('D20').Value = '0.0₃0679'
$ws.Range('E20').Value = '  -1.11%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.11'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.18'
$ws.Range('E23').Value = '  -1.85%  '
$ws.Range('E24').Value = '  +1.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.47'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.72'
$ws.Range('E26').Value = '  +1.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.89'
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('E29').Value = '  -1.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0466'
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('E31').Value = '  -4.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.15'
$ws.Range('E32').Value = '  -0.41%  '
$ws.Range('D33').Value = '1.403.95'
$ws.Range('E33').Value = '  +1.63%  '
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.914'
$ws.Range('E37').Value = '  -3.33%  '
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.526'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.809'
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.997'
$ws.Range('E42').Value = '  +0.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.45'
$ws.Range('E43').Value = '  +4.64%  '
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.18'
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.25'
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('D47').Value = '1.697.90'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.46'
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0503'
$ws.Range('E49').Value = '  +2.18%  '
$ws.Range('D50').Value = '0.0₇0979'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('E51').Value = '  +0.71%  '
